$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.081.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5028"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3897"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09211"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.88%  "
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.378"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.895.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.292"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.10%  "
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001108"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06640"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.224"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.139.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.323"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.106.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.541"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.077"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.93%  "
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.602"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.605"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.472"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06609"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.346"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02405"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2197"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.220"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6449"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.948"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6057"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.299"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.691"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("E51").Value = "  -1.44%  "
